$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    @(2609.2736383522001, 1805.9007177286148, 2046.4103443698491),
    @(2564.0016125361822, 1789.1326490189319, 2004.2483477758676),
    @(2772.704339868696, 2115.4570036857003, 1982.5372340769811),
    @(2536.9398168969196, 1807.1247017933108, 1671.066377746281),
    @(2848.1774558665747, 2142.9881146759785, 2201.7820875161801),
    @(2669.6642882999954, 2312.7943745309685, 2382.8080603810108),
    @(2404.1250818237031, 2125.2929879211752, 2073.5057216390678),
    @(2467.691911221747, 2069.4285611930673, 2059.5308696816064),
    @(3183.785045996166, 2098.4547370844157, 2032.8456309263427),
    @(2442.5675298805104, 1443.0734349368897, 1577.4517916709854),
    @(2285.2818322507096, 1678.0853514394139, 1621.0410838989533),
    @(2956.9778643979516, 2309.7218506473168, 2097.2723266668459),
    @(3174.1140970293991, 2482.204613906229, 2257.3095767361515),
    @(3188.7019285707479, 2538.6187625932635, 2288.5022984824991),
    @(3158.2049492718716, 2481.9163297845398, 2388.8584167805361),
    @(3056.1880566748437, 1719.0075552964661, 1829.9345956237353)
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i][0]
    $ws.Cells.Item($row, 2).Value = $values[$i][1]
    $ws.Cells.Item($row, 3).Value = $values[$i][2]
}
